$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 690.7059
$ws.Range("J17").Value = 702.63635
$ws.Range("L17").Value = 2107.90905
$ws.Range("N17").Value = -2443.90905
$ws.Range("H33").Value = 826.19354
$ws.Range("I33").Value = 300.5
$ws.Range("J33").Value = 5732.6665
$ws.Range("K33").Value = 300.5
$ws.Range("L33").Value = 5732.6665
$ws.Range("M33").Value = -71.5
$ws.Range("N33").Value = -6190.6665
$ws.Range("H43").Value = 3624.75
$ws.Range("I43").Value = 4499.6665
$ws.Range("J43").Value = 3099.8
$ws.Range("K43").Value = 4499.6665
$ws.Range("L43").Value = 3099.8
$ws.Range("M43").Value = -4430.6665
$ws.Range("N43").Value = -3237.8
$ws.Range("H74").Value = 3225
$ws.Range("I74").Value = 3225
$ws.Range("K74").Value = 3225
$ws.Range("M74").Value = -2289
$ws.Range("H77").Value = 3225
$ws.Range("I77").Value = 3225
$ws.Range("K77").Value = 16125
$ws.Range("M77").Value = -11445
$ws.Range("H129").Value = 3219.8445
$ws.Range("J129").Value = 1085.8649
$ws.Range("L129").Value = 3257.5947
$ws.Range("N129").Value = -13257.5947
$ws.Range("H132").Value = 2979236
$ws.Range("I132").Value = 3089470.5
$ws.Range("J132").Value = 2903.3333
$ws.Range("K132").Value = 9268411.5
$ws.Range("L132").Value = 8709.999899999999
$ws.Range("M132").Value = -9265881.5
$ws.Range("N132").Value = -13769.9999
$ws.Range("H136").Value = 52826.668
$ws.Range("J136").Value = 52826.668
$ws.Range("L136").Value = 52826.668
$ws.Range("N136").Value = -63026.668
$ws.Range("H138").Value = 6789.442
$ws.Range("I138").Value = 1317.2903
$ws.Range("K138").Value = 3951.8709
$ws.Range("M138").Value = 1188.1291

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30660.61
$ws.Range("I32").Value = 5080.96
$ws.Range("J32").Value = 172769.78
$ws.Range("K32").Value = 5080.96
$ws.Range("L32").Value = 172769.78
$ws.Range("M32").Value = -4793.96
$ws.Range("N32").Value = -173343.78
$ws.Range("H122").Value = 1548.4706
$ws.Range("I122").Value = 1524.7142
$ws.Range("J122").Value = 1659.3334
$ws.Range("K122").Value = 4574.142599999999
$ws.Range("L122").Value = 4978.0002
$ws.Range("M122").Value = -2124.142599999999
$ws.Range("N122").Value = -9878.0002
$ws.Range("H141").Value = 105000
$ws.Range("J141").Value = 105000
$ws.Range("L141").Value = 105000
$ws.Range("N141").Value = -115360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 13333.333
$ws.Range("I75").Value = 13333.333
$ws.Range("K75").Value = 13333.333
$ws.Range("M75").Value = -12397.333
$ws.Range("H78").Value = 13333.333
$ws.Range("I78").Value = 13333.333
$ws.Range("K78").Value = 39999.999
$ws.Range("M78").Value = -35319.999
$ws.Range("H134").Value = 2845.8604
$ws.Range("I134").Value = 2599.2974
$ws.Range("J134").Value = 4366.3335
$ws.Range("K134").Value = 7797.8922
$ws.Range("L134").Value = 13099.0005
$ws.Range("M134").Value = -5262.8922
$ws.Range("N134").Value = -18169.0005
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 35717730
$ws.Range("I132").Value = 31253288
$ws.Range("K132").Value = 93759864
$ws.Range("M132").Value = -93757334
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 843.6900000000001
$ws.Range("I131").Value = 373.33334
$ws.Range("J131").Value = 873.71277
$ws.Range("K131").Value = 1120.00002
$ws.Range("L131").Value = 2621.13831
$ws.Range("M131").Value = 3919.99998
$ws.Range("N131").Value = -12701.13831
$ws.Range("H132").Value = 2074.6875
$ws.Range("J132").Value = 2617.7273
$ws.Range("L132").Value = 23559.5457
$ws.Range("N132").Value = -28619.5457
$ws.Range("H139").Value = 2926.25
$ws.Range("I139").Value = 2728.5715
$ws.Range("J139").Value = 3007.647
$ws.Range("K139").Value = 8185.7145
$ws.Range("L139").Value = 9022.940999999999
$ws.Range("M139").Value = -3045.7145
$ws.Range("N139").Value = -19302.941

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 71107.8
$ws.Range("I70").Value = 128936.875
$ws.Range("J70").Value = 5017.4287
$ws.Range("K70").Value = 128936.875
$ws.Range("L70").Value = 5017.4287
$ws.Range("M70").Value = -128666.875
$ws.Range("N70").Value = -5557.4287
$ws.Range("H73").Value = 71107.8
$ws.Range("I73").Value = 128936.875
$ws.Range("J73").Value = 5017.4287
$ws.Range("K73").Value = 128936.875
$ws.Range("L73").Value = 5017.4287
$ws.Range("M73").Value = -128000.875
$ws.Range("N73").Value = -6889.4287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6095.923
$ws.Range("I132").Value = 6113.5454
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 18340.6362
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -15810.6362
$ws.Range("N132").Value = -23057
$ws.Range("H136").Value = 1650.4546
$ws.Range("I136").Value = 1283.0769
$ws.Range("K136").Value = 3849.2307
$ws.Range("M136").Value = -1299.2307
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4170430
$ws.Range("I62").Value = 31251250
$ws.Range("J62").Value = 4150
$ws.Range("K62").Value = 31251250
$ws.Range("L62").Value = 4150
$ws.Range("M62").Value = -31250626
$ws.Range("N62").Value = -5398
$ws.Range("H65").Value = 4170430
$ws.Range("I65").Value = 31251250
$ws.Range("J65").Value = 4150
$ws.Range("K65").Value = 156256250
$ws.Range("L65").Value = 20750
$ws.Range("M65").Value = -156253130
$ws.Range("N65").Value = -26990
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26622
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -83112
$ws.Range("H132").Value = 3545.6538
$ws.Range("I132").Value = 4331.8125
$ws.Range("J132").Value = 2287.8
$ws.Range("K132").Value = 12995.4375
$ws.Range("L132").Value = 6863.400000000001
$ws.Range("M132").Value = -10465.4375
$ws.Range("N132").Value = -11923.4
$ws.Range("H136").Value = 987.5
$ws.Range("I136").Value = 588.5
$ws.Range("J136").Value = 2317.5
$ws.Range("K136").Value = 1765.5
$ws.Range("L136").Value = 6952.5
$ws.Range("M136").Value = 784.5
$ws.Range("N136").Value = -12052.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
